# "case 3/storage&modifiers.xlsx" update
# Test case data refresh: test case #1 re-numbered/re-worded, a second
# full test case (#2) appended below it, column widths widened, the
# F2:F13 "Result" merge replaced by two standalone PASS cells (F8/F28),
# and the old centered look of column F dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The old "Result" column was one big merged/centered cell (F2:F13).
# Split it apart and drop the centering so each row stands on its own.
# ---------------------------------------------------------------------
$ws.Range("F2:F13").UnMerge()
$ws.Range("F2:F15").HorizontalAlignment = 1

# ---------------------------------------------------------------------
# Widen columns D and E a bit.
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 35.166666666666664
$ws.Columns("E").ColumnWidth = 35.666666666666664

# ---------------------------------------------------------------------
# Re-number test case 1 and refresh its Given/Expected/Actual values.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("D2").Value = 'welcome to hcl'
$ws.Range("E2").Value = 'welcome to hcl'
$ws.Range("F2").Value = ''

$ws.Range("D3").Value = 'static=1'
$ws.Range("E3").Value = 'static=1'

$ws.Range("D4").Value = 'non-static=1'
$ws.Range("E4").Value = 'non-static=1'

$ws.Range("D5").Value = 'register=5 size=4'
$ws.Range("E5").Value = 'register=5 size=4'

$ws.Range("C6").Value = 'short int=2'
$ws.Range("D6").Value = 'double=2.3 size=16'
$ws.Range("E6").Value = 'double=2.3 size=16'

$ws.Range("C7").Value = 'extern iX=10'
$ws.Range("D7").Value = 'global=96 size=4'
$ws.Range("E7").Value = 'global=96 size=4'

$ws.Range("D8").Value = 'second global=100 size=4'
$ws.Range("E8").Value = 'second global=100 size=4'
$ws.Range("F8").Value = 'PASS'

$ws.Range("D9").Value = 'global variables in function 12,-96'
$ws.Range("E9").Value = 'global variables in function 12,-96'

$ws.Range("D10").Value = 'global variables outside function 96,100'
$ws.Range("E10").Value = 'global variables outside function 96,100'

$ws.Range("D11").Value = 'static=2'
$ws.Range("E11").Value = 'static=2'

$ws.Range("D12").Value = 'non static=1'
$ws.Range("E12").Value = 'non static=1'

$ws.Range("D13").Value = 'iX=10'
$ws.Range("E13").Value = 'iX=10'

# New rows 14-15 continue test case 1.
$ws.Range("D14").Value = 'static=3'
$ws.Range("E14").Value = 'static=3'

$ws.Range("D15").Value = 'non static=1'
$ws.Range("E15").Value = 'non static=1'

# ---------------------------------------------------------------------
# New test case 2 block (rows 20-33), same structure as case 1.
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 'double=5.3'
$ws.Range("D20").Value = 'welcome to hcl'
$ws.Range("E20").Value = 'welcome to hcl'

$ws.Range("C21").Value = 'register=15'
$ws.Range("D21").Value = 'static=1'
$ws.Range("E21").Value = 'static=1'

$ws.Range("A22").Value = 2
$ws.Range("B22").Value = 'with neg value for unsigned variable'
$ws.Range("C22").Value = 'unsigned=-9687'
$ws.Range("D22").Value = 'non-static=1'
$ws.Range("E22").Value = 'non-static=1'

$ws.Range("C23").Value = 'signed=-96'
$ws.Range("D23").Value = 'register=5 size=4'
$ws.Range("E23").Value = 'register=5 size=4'

$ws.Range("C24").Value = 'short int=2'
$ws.Range("D24").Value = 'double=2.3 size=16'
$ws.Range("E24").Value = 'double=2.3 size=16'

$ws.Range("C25").Value = 'extern iX=10'
$ws.Range("D25").Value = 'global=96 size=4'
$ws.Range("E25").Value = 'global=96 size=4'

$ws.Range("D26").Value = 'second global=100 size=4'
$ws.Range("E26").Value = 'second global=100 size=4'

$ws.Range("D27").Value = 'global variables in function garbage,-96'
$ws.Range("E27").Value = 'global variables in function garbage,-96'

$ws.Range("D28").Value = 'global variables outside function 96,100'
$ws.Range("E28").Value = 'global variables outside function 96,100'
$ws.Range("F28").Value = 'PASS'

$ws.Range("D29").Value = 'static=2'
$ws.Range("E29").Value = 'static=2'

$ws.Range("D30").Value = 'non static=1'
$ws.Range("E30").Value = 'non static=1'

$ws.Range("D31").Value = 'iX=10'
$ws.Range("E31").Value = 'iX=10'

$ws.Range("D32").Value = 'static=3'
$ws.Range("E32").Value = 'static=3'

$ws.Range("D33").Value = 'non static=1'
$ws.Range("E33").Value = 'non static=1'

# ---------------------------------------------------------------------
# Leave the selection on the new PASS cell, scrolled down to it.
# ---------------------------------------------------------------------
$ws.Range("F28").Select()
